# Ajout des taux d'accroissement (Population et PIB par habitant)
# P = Taux accroiss Population = (B[r]-B[r-1])/B[r-1]*100
# Q = Taux accroiss PIB        = (N[r]-N[r-1])/N[r-1]*100

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers (P1, Q1) : same style (s="1") as the other header cells ---
$ws.Range("O1").Copy($ws.Range("P1"))
$ws.Range("P1").Value = "Taux accroiss Population"

$ws.Range("O1").Copy($ws.Range("Q1"))
$ws.Range("Q1").Value = "Taux accroiss PIB"

# Row 2 (1962) has no prior year, so P2/Q2 stay blank (no prior-year growth rate).

# --- Data rows 3..63 ---
$data = @(
    @(3, 1.852411390512154, 9.970958099789119),
    @(4, 1.925999695428571, 5.181413503482668),
    @(5, 1.994561396566508, 2.987860694167144),
    @(6, 2.044318649888899, 1.916501943455562),
    @(7, 2.085052963116807, -0.5354274781566137),
    @(8, 2.132268150645489, 0.8792861413132735),
    @(9, 2.17993856707166, 0.0253902930783978),
    @(10, 2.229154478887807, -0.06933550214139883),
    @(11, 2.274345958329649, -0.8466397374539025),
    @(12, 2.321576053339869, 6.554044594604358),
    @(13, 2.367949548922654, 0.1756615742589895),
    @(14, 2.410144553244464, -9.631429332847564),
    @(15, 2.461582932184392, 1.527658133522558),
    @(16, 2.518976713055965, 1.948080621982795),
    @(17, 2.576470841368939, 2.904257369893015),
    @(18, 2.624832951626765, 0.291492466254617),
    @(19, 2.695394810269724, 1.801502139996458),
    @(20, 3.196780246327857, -0.7361991832215198),
    @(21, 2.899917692680876, -0.5763285310710353),
    @(22, 2.938482432553169, -0.5749280374805377),
    @(23, 2.973254918159052, -0.5548149536395264),
    @(24, 3.0043796952264, -0.517752003162586),
    @(25, 3.032007665413428, -0.4654866066174934),
    @(26, 3.056294153665351, -0.3997566988461876),
    @(27, 3.077397191895126, -0.3222906799551839),
    @(28, 3.095475957045557, -0.2348019334505569),
    @(29, 3.110689418905555, -0.1389793034266895),
    @(30, 3.123195172000659, -0.03647453289306046),
    @(31, 3.133148425993437, 0.07111228798437619),
    @(32, 3.140701160893289, 0.182248232355664),
    @(33, 3.146001437827861, 0.2954826183072301),
    @(34, 3.149192836566606, 0.4094615340172858),
    @(35, 3.150414015172776, 0.522940078818146),
    @(36, 3.149798380391999, 0.6347923891886653),
    @(37, 3.147473847524251, 0.7440189310775391),
    @(38, 3.143562697162627, 0.8497505936860605),
    @(39, 2.906151181259209, 3.733039025005702),
    @(40, 3.12447074877571, 2.838707853544364),
    @(41, 3.101479045072852, 1.999548844771426),
    @(42, 3.086791171586545, 1.51752867785675),
    @(43, 3.096166909823639, 0.3497930373576263),
    @(44, 3.109084084937952, 1.382537855867638),
    @(45, 3.273034971775046, -1.621780028301745),
    @(46, 3.169872245240524, 0.8481851881178004),
    @(47, 2.99919178598731, 2.959561487677131),
    @(48, 3.065646083876339, 1.768132925069521),
    @(49, 3.04446047696334, -0.6932574407960579),
    @(50, 3.029216388838729, -0.9142356841476795),
    @(51, 3.022602537549424, -0.04719822755828718),
    @(52, 3.012139446856699, 1.631847815600751),
    @(53, 2.992760772931846, 4.125162298328844),
    @(54, 2.999112695182315, 3.331830131562641),
    @(55, 2.997977153952069, -1.297448523688027),
    @(56, 2.967788638726843, 0.4367820263125211),
    @(57, 2.923086155549037, 2.604237999880965),
    @(58, 2.853598653795997, 3.707277200582371),
    @(59, 2.773143642021814, 4.003586262137659),
    @(60, 2.698362622679551, 1.207015223707297),
    @(61, 2.626194045386865, 4.377109446420091),
    @(62, 2.580133011595787, 3.5449095030474),
    @(63, 2.55483828955716, 3.31407470119367)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 16).Value = $row[1]
    $ws.Cells.Item($r, 17).Value = $row[2]
}

Write-Output "Added Taux accroiss Population (P) and Taux accroiss PIB (Q) columns for rows 3-63"
